$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-unneeded rows 7-36 (each planeswalker used to span 7 rows,
# now each planeswalker is consolidated into a single row)
$ws.Rows("7:36").Delete()

$ws.Range("A2").Value = '(''Ajani, Caller of the Pride'', [''{1}{W}{W}'', ''Legendary Planeswalker — Ajani'', ''+1: Put a +1/+1 counter on up to one target creature.'', ''−3: Target creature gains flying and double strike until end of turn.'', ''−8: Create X 2/2 white Cat creature tokens, where X is your life total.'', ''Loyalty: 4''])'
$ws.Range("A3").Value = '(''Chandra, Pyromaster'', [''{2}{R}{R}'', ''Legendary Planeswalker — Chandra'', ''+1: Chandra, Pyromaster deals 1 damage to target player or planeswalker and 1 damage to up to one target creature that player or that planeswalker’s controller controls. That creature can’t block this turn.'', ''0: Exile the top card of your library. You may play it this turn.'', ''−7: Exile the top ten cards of your library. Choose an instant or sorcery card exiled this way and copy it three times. You may cast the copies without paying their mana costs.'', ''Loyalty: 4''])'
$ws.Range("A4").Value = '(''Garruk, Caller of Beasts'', [''{4}{G}{G}'', ''Legendary Planeswalker — Garruk'', ''+1: Reveal the top five cards of your library. Put all creature cards revealed this way into your hand and the rest on the bottom of your library in any order.'', ''−3: You may put a green creature card from your hand onto the battlefield.'', ''−7: You get an emblem with “Whenever you cast a creature spell, you may search your library for a creature card, put it onto the battlefield, then shuffle your library.”'', ''Loyalty: 4''])'
$ws.Range("A5").Value = '(''Jace, Memory Adept'', [''{3}{U}{U}'', ''Legendary Planeswalker — Jace'', ''+1: Draw a card. Target player mills a card.'', ''0: Target player mills ten cards.'', ''−7: Any number of target players each draw twenty cards.'', ''Loyalty: 4''])'
$ws.Range("A6").Value = '(''Liliana of the Dark Realms'', [''{2}{B}{B}'', ''Legendary Planeswalker — Liliana'', ''+1: Search your library for a Swamp card, reveal it, and put it into your hand. Then shuffle your library.'', ''−3: Target creature gets +X/+X or -X/-X until end of turn, where X is the number of Swamps you control.'', ''−6: You get an emblem with “Swamps you control have ‘{T}: Add {B}{B}{B}{B}.’”'', ''Loyalty: 3''])'

